$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.003455333333333333
$ws.Range("H2").Value = 0.010366
$ws.Range("I2").Value = 0.000270121469710956
$ws.Range("J2").Value = 0.000270121469710956
$ws.Range("M2").Value = 3.618510333333333
$ws.Range("N2").Value = 10.855531
$ws.Range("O2").Value = 0.1815566256530994
$ws.Range("P2").Value = 0.1815566256530994
$ws.Range("Q2").Value = 0.01250315937177778
$ws.Range("R2").Value = 0.112528434346
$ws.Range("S2").Value = 0.00004904234255717705
$ws.Range("T2").Value = 0.00004904234255717707
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.003455333333333333
$ws.Range("H3").Value = 0.010366
$ws.Range("I3").Value = 0.000270121469710956
$ws.Range("J3").Value = 0.000270121469710956
$ws.Range("O3").Value = 0.1937079481987336
$ws.Range("P3").Value = 0.1937079481987336
$ws.Range("Q3").Value = 0.01333997775733333
$ws.Range("R3").Value = 0.120059799816
$ws.Range("S3").Value = 0.00005232467566213565
$ws.Range("T3").Value = 0.00005232467566213567
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.003455333333333333
$ws.Range("H4").Value = 0.010366
$ws.Range("I4").Value = 0.000270121469710956
$ws.Range("J4").Value = 0.000270121469710956
$ws.Range("M4").Value = 12.45127566666667
$ws.Range("N4").Value = 37.353827
$ws.Range("O4").Value = 0.6247354261481669
$ws.Range("P4").Value = 0.6247354261481669
$ws.Range("Q4").Value = 0.04302330785355555
$ws.Range("R4").Value = 0.387209770682
$ws.Range("S4").Value = 0.0001687544514916432
$ws.Range("T4").Value = 0.0001687544514916433
$ws.Range("I5").Value = 0.8991926531546518
$ws.Range("J5").Value = 0.8991926531546519
$ws.Range("M5").Value = 3.618510333333333
$ws.Range("N5").Value = 10.855531
$ws.Range("O5").Value = 0.1815566256530994
$ws.Range("P5").Value = 0.1815566256530994
$ws.Range("Q5").Value = 41.62108647030033
$ws.Range("R5").Value = 374.5897782327029
$ws.Range("S5").Value = 0.1632543839188164
$ws.Range("T5").Value = 0.1632543839188164
$ws.Range("I6").Value = 0.8991926531546518
$ws.Range("J6").Value = 0.8991926531546519
$ws.Range("O6").Value = 0.1937079481987336
$ws.Range("P6").Value = 0.1937079481987336
$ws.Range("S6").Value = 0.1741807638779631
$ws.Range("T6").Value = 0.1741807638779632
$ws.Range("I7").Value = 0.8991926531546518
$ws.Range("J7").Value = 0.8991926531546519
$ws.Range("M7").Value = 12.45127566666667
$ws.Range("N7").Value = 37.353827
$ws.Range("O7").Value = 0.6247354261481669
$ws.Range("P7").Value = 0.6247354261481669
$ws.Range("Q7").Value = 143.2179470137056
$ws.Range("R7").Value = 1288.961523123351
$ws.Range("S7").Value = 0.5617575053578723
$ws.Range("T7").Value = 0.5617575053578723
$ws.Range("G8").Value = 1.286049666666667
$ws.Range("H8").Value = 3.858149
$ws.Range("I8").Value = 0.1005372253756372
$ws.Range("J8").Value = 0.1005372253756372
$ws.Range("M8").Value = 3.618510333333333
$ws.Range("N8").Value = 10.855531
$ws.Range("O8").Value = 0.1815566256530994
$ws.Range("P8").Value = 0.1815566256530994
$ws.Range("Q8").Value = 4.653584008013222
$ws.Range("R8").Value = 41.882256072119
$ws.Range("S8").Value = 0.01825319939172585
$ws.Range("T8").Value = 0.01825319939172585
$ws.Range("G9").Value = 1.286049666666667
$ws.Range("H9").Value = 3.858149
$ws.Range("I9").Value = 0.1005372253756372
$ws.Range("J9").Value = 0.1005372253756372
$ws.Range("O9").Value = 0.1937079481987336
$ws.Range("P9").Value = 0.1937079481987336
$ws.Range("Q9").Value = 4.965041659702667
$ws.Range("R9").Value = 44.68537493732401
$ws.Range("S9").Value = 0.01947485964510834
$ws.Range("T9").Value = 0.01947485964510834
$ws.Range("G10").Value = 1.286049666666667
$ws.Range("H10").Value = 3.858149
$ws.Range("I10").Value = 0.1005372253756372
$ws.Range("J10").Value = 0.1005372253756372
$ws.Range("M10").Value = 12.45127566666667
$ws.Range("N10").Value = 37.353827
$ws.Range("O10").Value = 0.6247354261481669
$ws.Range("P10").Value = 0.6247354261481669
$ws.Range("Q10").Value = 16.01295892069145
$ws.Range("R10").Value = 144.116630286223
$ws.Range("S10").Value = 0.062809166338803
$ws.Range("T10").Value = 0.06280916633880301
